# GPLIM-2588 Fix spreadsheet headers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the header cells with corrected labels.
# A1 originally carried a quote-prefix (forced-text) style, so re-enter its
# replacement text the same way to avoid Excel allocating a new cell style.
$ws.Range("A1").Value = "'Specimen_Number"
$ws.Range("F1").Value = "SAMPLE_TYPE"

# Leave F2 selected, matching the saved selection state.
$ws.Range("F2").Select()
